$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Cells.Item(28, 8).Value = 2431.1667
$ws.Cells.Item(28, 9).Value = 1478.625
$ws.Cells.Item(28, 10).Value = 3193.2
$ws.Cells.Item(28, 11).Value = 1478.625
$ws.Cells.Item(28, 12).Value = 3193.2
$ws.Cells.Item(28, 13).Value = -993.625
$ws.Cells.Item(28, 14).Value = -4163.2
# row 33
$ws.Cells.Item(33, 8).Value = 10888.429
$ws.Cells.Item(33, 10).Value = 1683.8
$ws.Cells.Item(33, 12).Value = 1683.8
$ws.Cells.Item(33, 14).Value = -2141.8
# row 39
$ws.Cells.Item(39, 8).Value = 1328.5385
$ws.Cells.Item(39, 10).Value = 4250
$ws.Cells.Item(39, 12).Value = 12750
$ws.Cells.Item(39, 14).Value = -13342
# row 113
$ws.Cells.Item(113, 8).Value = 7192.25
$ws.Cells.Item(113, 9).Value = 6250
$ws.Cells.Item(113, 10).Value = 7757.6
$ws.Cells.Item(113, 11).Value = 6250
$ws.Cells.Item(113, 12).Value = 7757.6
$ws.Cells.Item(113, 13).Value = -2996
$ws.Cells.Item(113, 14).Value = -14265.6
# row 127
$ws.Cells.Item(127, 8).Value = 783.1818
$ws.Cells.Item(127, 9).Value = 685
$ws.Cells.Item(127, 10).Value = 1225
$ws.Cells.Item(127, 11).Value = 2055
$ws.Cells.Item(127, 12).Value = 3675
$ws.Cells.Item(127, 13).Value = 2905
$ws.Cells.Item(127, 14).Value = -13595
# row 138
$ws.Cells.Item(138, 8).Value = 25642764
$ws.Cells.Item(138, 9).Value = 1186.2858
$ws.Cells.Item(138, 10).Value = 90912230
$ws.Cells.Item(138, 11).Value = 3558.8574
$ws.Cells.Item(138, 12).Value = 272736690
$ws.Cells.Item(138, 13).Value = 1581.1426
$ws.Cells.Item(138, 14).Value = -272746970

$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Cells.Item(61, 8).Value = 2415.1
$ws.Cells.Item(61, 9).Value = 2111.0386
$ws.Cells.Item(61, 10).Value = 4391.5
$ws.Cells.Item(61, 11).Value = 2111.0386
$ws.Cells.Item(61, 12).Value = 4391.5
$ws.Cells.Item(61, 13).Value = -1899.0386
$ws.Cells.Item(61, 14).Value = -4815.5
# row 102
$ws.Cells.Item(102, 8).Value = 6070.7334
$ws.Cells.Item(102, 9).Value = 5556.091
$ws.Cells.Item(102, 11).Value = 5556.091
$ws.Cells.Item(102, 13).Value = -3934.091
# row 132
$ws.Cells.Item(132, 8).Value = 31888.447
$ws.Cells.Item(132, 9).Value = 2223.439
$ws.Cells.Item(132, 11).Value = 6670.316999999999
$ws.Cells.Item(132, 13).Value = -4140.316999999999
# row 136
$ws.Cells.Item(136, 8).Value = 2415.1
$ws.Cells.Item(136, 9).Value = 2111.0386
$ws.Cells.Item(136, 10).Value = 4391.5
$ws.Cells.Item(136, 11).Value = 6333.1158
$ws.Cells.Item(136, 12).Value = 13174.5
$ws.Cells.Item(136, 13).Value = -3783.1158
$ws.Cells.Item(136, 14).Value = -18274.5

$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Cells.Item(99, 8).Value = 3758.5557
$ws.Cells.Item(99, 9).Value = 3215
$ws.Cells.Item(99, 11).Value = 3215
$ws.Cells.Item(99, 13).Value = -1717
# row 113
$ws.Cells.Item(113, 8).Value = 8166.5
$ws.Cells.Item(113, 9).Value = 8166.5
$ws.Cells.Item(113, 11).Value = 8166.5
$ws.Cells.Item(113, 13).Value = -5996.5
# row 141
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Cells.Item(22, 8).Value = 6650.875
$ws.Cells.Item(22, 9).Value = 7532.143
$ws.Cells.Item(22, 10).Value = 482
$ws.Cells.Item(22, 11).Value = 7532.143
$ws.Cells.Item(22, 12).Value = 482
$ws.Cells.Item(22, 13).Value = -7182.143
$ws.Cells.Item(22, 14).Value = -1182
# row 31
$ws.Cells.Item(31, 8).Value = 3896.1836
$ws.Cells.Item(31, 9).Value = 2504.88
$ws.Cells.Item(31, 10).Value = 5345.4585
$ws.Cells.Item(31, 11).Value = 2504.88
$ws.Cells.Item(31, 12).Value = 5345.4585
$ws.Cells.Item(31, 13).Value = -2209.88
$ws.Cells.Item(31, 14).Value = -5935.4585
# row 34
$ws.Cells.Item(34, 8).Value = 3896.1836
$ws.Cells.Item(34, 9).Value = 2504.88
$ws.Cells.Item(34, 10).Value = 5345.4585
$ws.Cells.Item(34, 11).Value = 2504.88
$ws.Cells.Item(34, 12).Value = 5345.4585
$ws.Cells.Item(34, 13).Value = -2302.88
$ws.Cells.Item(34, 14).Value = -5749.4585
# row 99
$ws.Cells.Item(99, 8).Value = 7312.375
$ws.Cells.Item(99, 9).Value = 7357
$ws.Cells.Item(99, 11).Value = 7357
$ws.Cells.Item(99, 13).Value = -5859
# row 126
$ws.Cells.Item(126, 8).Value = 7312.375
$ws.Cells.Item(126, 9).Value = 7357
$ws.Cells.Item(126, 11).Value = 22071
$ws.Cells.Item(126, 13).Value = -19601
# row 131
$ws.Cells.Item(131, 8).Value = 28052.875
$ws.Cells.Item(131, 10).Value = 30631.857
$ws.Cells.Item(131, 12).Value = 30631.857
$ws.Cells.Item(131, 14).Value = -40711.857
# row 132
$ws.Cells.Item(132, 8).Value = 4495.879
$ws.Cells.Item(132, 9).Value = 4344.0835
$ws.Cells.Item(132, 11).Value = 13032.2505
$ws.Cells.Item(132, 13).Value = -10502.2505

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Cells.Item(5, 8).Value = 2537.8
$ws.Cells.Item(5, 10).Value = 7747.5
$ws.Cells.Item(5, 12).Value = 23242.5
$ws.Cells.Item(5, 14).Value = -23466.5
# row 12
$ws.Cells.Item(12, 8).Value = 770
$ws.Cells.Item(12, 9).Value = 548.5
$ws.Cells.Item(12, 11).Value = 1645.5
$ws.Cells.Item(12, 13).Value = -1472.5
# row 41
$ws.Cells.Item(41, 8).Value = 1549.5
$ws.Cells.Item(41, 10).Value = 1460
$ws.Cells.Item(41, 12).Value = 4380
$ws.Cells.Item(41, 14).Value = -5056
# row 68
$ws.Cells.Item(68, 8).Value = 633.6
$ws.Cells.Item(68, 10).Value = 624.5
$ws.Cells.Item(68, 12).Value = 1873.5
$ws.Cells.Item(68, 14).Value = -3495.5
# row 71
$ws.Cells.Item(71, 8).Value = 633.6
$ws.Cells.Item(71, 10).Value = 624.5
$ws.Cells.Item(71, 12).Value = 5620.5
$ws.Cells.Item(71, 14).Value = -13732.5
# row 82
$ws.Cells.Item(82, 8).Value = 3913
$ws.Cells.Item(82, 9).Value = 3913
$ws.Cells.Item(82, 11).Value = 11739
$ws.Cells.Item(82, 13).Value = -11333
# row 85
$ws.Cells.Item(85, 8).Value = 3913
$ws.Cells.Item(85, 9).Value = 3913
$ws.Cells.Item(85, 11).Value = 11739
$ws.Cells.Item(85, 13).Value = -10335
# row 87
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 13).ClearContents()
# row 90
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 13).ClearContents()
# row 104
$ws.Cells.Item(104, 8).Value = 2900
$ws.Cells.Item(104, 10).Value = 2900
$ws.Cells.Item(104, 12).Value = 8700
$ws.Cells.Item(104, 14).Value = -13942
# row 135
$ws.Cells.Item(135, 8).Value = 2537.8
$ws.Cells.Item(135, 10).Value = 7747.5
$ws.Cells.Item(135, 12).Value = 69727.5
$ws.Cells.Item(135, 14).Value = -74797.5

$ws = $wb.Worksheets.Item("GSM")
# row 14
$ws.Cells.Item(14, 8).Value = 11762500
$ws.Cells.Item(14, 9).Value = 11762500
$ws.Cells.Item(14, 11).Value = 11762500
$ws.Cells.Item(14, 13).Value = -11762332
# row 18
$ws.Cells.Item(18, 8).Value = 2512500
$ws.Cells.Item(18, 9).Value = 2512500
$ws.Cells.Item(18, 11).Value = 2512500
$ws.Cells.Item(18, 13).Value = -2512207
# row 44
$ws.Cells.Item(44, 8).Value = 21999.2
$ws.Cells.Item(44, 9).Value = 17499
$ws.Cells.Item(44, 10).Value = 24999.334
$ws.Cells.Item(44, 11).Value = 17499
$ws.Cells.Item(44, 12).Value = 24999.334
$ws.Cells.Item(44, 13).Value = -16903
$ws.Cells.Item(44, 14).Value = -26191.334
# row 109
$ws.Cells.Item(109, 8).Value = 40747.5
$ws.Cells.Item(109, 9).Value = 34500
$ws.Cells.Item(109, 10).Value = 42830
$ws.Cells.Item(109, 11).Value = 34500
$ws.Cells.Item(109, 12).Value = 42830
$ws.Cells.Item(109, 13).Value = -33460
$ws.Cells.Item(109, 14).Value = -44910
# row 113
$ws.Cells.Item(113, 8).Value = 3500
$ws.Cells.Item(113, 9).Value = 3500
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 3500
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -1330
$ws.Cells.Item(113, 14).ClearContents()
# row 126
$ws.Cells.Item(126, 8).Value = 501000
$ws.Cells.Item(126, 9).Value = 1000000
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 3000000
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -2997530
$ws.Cells.Item(126, 14).Value = -10940
# row 132
$ws.Cells.Item(132, 8).Value = 1207
$ws.Cells.Item(132, 9).Value = 1187.5
$ws.Cells.Item(132, 10).Value = 1233
$ws.Cells.Item(132, 11).Value = 3562.5
$ws.Cells.Item(132, 12).Value = 3699
$ws.Cells.Item(132, 13).Value = -1032.5
$ws.Cells.Item(132, 14).Value = -8759

$ws = $wb.Worksheets.Item("LTW")
# row 43
$ws.Cells.Item(43, 8).Value = 23141.428
$ws.Cells.Item(43, 10).Value = 22998.334
$ws.Cells.Item(43, 12).Value = 22998.334
$ws.Cells.Item(43, 14).Value = -23384.334
# row 122
$ws.Cells.Item(122, 8).Value = 4593.2
$ws.Cells.Item(122, 9).Value = 4488.6665
$ws.Cells.Item(122, 10).Value = 4750
$ws.Cells.Item(122, 11).Value = 13465.9995
$ws.Cells.Item(122, 12).Value = 14250
$ws.Cells.Item(122, 13).Value = -11015.9995
$ws.Cells.Item(122, 14).Value = -19150
# row 132
$ws.Cells.Item(132, 8).Value = 15550.071
$ws.Cells.Item(132, 9).Value = 2418.889
$ws.Cells.Item(132, 10).Value = 39186.2
$ws.Cells.Item(132, 11).Value = 7256.667
$ws.Cells.Item(132, 12).Value = 117558.6
$ws.Cells.Item(132, 13).Value = -4726.667
$ws.Cells.Item(132, 14).Value = -122618.6
# row 136
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# row 136
$ws.Cells.Item(136, 8).Value = 5109.7334
$ws.Cells.Item(136, 9).Value = 4175.2
$ws.Cells.Item(136, 10).Value = 6978.8
$ws.Cells.Item(136, 11).Value = 12525.6
$ws.Cells.Item(136, 12).Value = 20936.4
$ws.Cells.Item(136, 13).Value = -9975.599999999999
$ws.Cells.Item(136, 14).Value = -26036.4

